$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the other workers' account-statement rows (22-27): this
#    deletes the 6 entire rows, shifting the former "last row" (28,
#    with its heavier bottom-border style) up to row 22, and the
#    signature block (formerly rows 33/34) up to rows 27/28.
# ------------------------------------------------------------------
$ws.Range("B22:J27").EntireRow.Delete()

# ------------------------------------------------------------------
# 2. Update the header summary figures.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 258046   # VALOR MORA
$ws.Range("C13").Value = 1        # Cant. Trabajadores
$ws.Range("F13").Value = 7        # Cant. Periodos

# ------------------------------------------------------------------
# 3. Rewrite the account-statement table body (rows 16-22) so it only
#    lists the single remaining worker, periods sorted descending
#    (2112 .. 2106), matching the new data export.
# ------------------------------------------------------------------
$periods = @(2112, 2111, 2110, 2109, 2108, 2107, 2106)
$valores = @(40000, 36341, 36341, 36341, 36341, 36341, 36341)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "1002244410"
    $ws.Range("D$r").Value = "JESUS DAVID GONZALEZ MARRUGO"
    $ws.Range("E$r").Value = [string]$periods[$i]
    $ws.Range("F$r").Value = $valores[$i]
    $ws.Range("G$r").Value = 908526
}
